$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the "Run No." column (B): "Norun" -> "NoRun" for every
# test-case row (rows 2-15) on the OPD billing / smoke-sanity test script.
for ($r = 2; $r -le 15; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value() -eq "Norun") {
        $cell.Value = "NoRun"
    }
}

# Leave the sheet with the same cursor/selection state as the saved file.
[void]$ws.Range("B20").Select()
